{"js": "// Update \"resultados\" table:\n//  - n changes from 365 to 364\n//  - second header column label \"Grupo (%)\" -> \"CAE (%)\"\n//  - the \"CIRURGIA VASCULAR\" row is removed entirely\n//  - a few percentages are corrected (rounding fixes)\n\nasync function replaceText(oldText, newText) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1) n: 365 -> 364\nawait replaceText(\"365\", \"364\");\n\n// 2) Header label: \"Grupo (%)\" -> \"CAE (%)\"\nawait replaceText(\"Grupo (%)\", \"CAE (%)\");\n\n// 3) Remove the \"CIRURGIA VASCULAR\" row entirely.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.rows.load(\"items\");\ntable.load(\"values\");\nawait context.sync();\n\nlet targetRowIndex = -1;\nfor (let i = 0; i < table.values.length; i++) {\n  if (table.values[i][0] === \"CIRURGIA VASCULAR\") {\n    targetRowIndex = i;\n    break;\n  }\n}\nif (targetRowIndex === -1) {\n  throw new Error(\"Could not find the CIRURGIA VASCULAR row\");\n}\ntable.rows.items[targetRowIndex].delete();\nawait context.sync();\n\n// 4) Percentage corrections (rounding fixes)\nawait replaceText(\"134 (36.7)\", \"134 (36.8)\");\nawait replaceText(\"84 (23.0)\", \"84 (23.1)\");\nawait replaceText(\"89 (24.4)\", \"89 (24.5)\");\n", "ps1": "# Update \"resultados\" table:\n#  - n changes from 365 to 364\n#  - second header column label \"Grupo (%)\" -> \"CAE (%)\"\n#  - the \"CIRURGIA VASCULAR\" row is removed entirely\n#  - a few percentages are corrected (rounding fixes)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nfunction Clean-Text($text) {\n  return $text.Replace([char]13, \"\").Replace([char]7, \"\")\n}\n\n# Set the text of the 2nd-column cell of the row whose 1st-column label matches $label.\nfunction Set-RowValue($table, $label, $newValue) {\n  for ($i = 1; $i -le $table.Rows.Count; $i++) {\n    $cellText = Clean-Text $table.Cell($i, 1).Range.Text\n    if ($cellText -eq $label) {\n      $table.Cell($i, 2).Range.Text = $newValue\n      return\n    }\n  }\n  throw \"Row not found: $label\"\n}\n\n# Set the text of the 1st-column cell of the row whose 1st-column label matches $label.\nfunction Set-RowLabel($table, $label, $newLabel) {\n  for ($i = 1; $i -le $table.Rows.Count; $i++) {\n    $cellText = Clean-Text $table.Cell($i, 1).Range.Text\n    if ($cellText -eq $label) {\n      $table.Cell($i, 1).Range.Text = $newLabel\n      return\n    }\n  }\n  throw \"Row not found: $label\"\n}\n\n# 1) n: 365 -> 364\nSet-RowValue $t \"n\" \"364\"\n\n# 2) Header label: \"Grupo (%)\" -> \"CAE (%)\"\nSet-RowLabel $t \"Grupo (%)\" \"CAE (%)\"\n\n# 3) Remove the \"CIRURGIA VASCULAR\" row entirely (look it up dynamically).\n$removed = $false\nfor ($i = $t.Rows.Count; $i -ge 1; $i--) {\n  $cellText = Clean-Text $t.Cell($i, 1).Range.Text\n  if ($cellText -eq \"CIRURGIA VASCULAR\") {\n    $t.Rows.Item($i).Delete()\n    $removed = $true\n    break\n  }\n}\nif (-not $removed) {\n  throw \"Could not find the CIRURGIA VASCULAR row\"\n}\n\n# 4) Percentage corrections (rounding fixes)\nSet-RowValue $t \"JOELHO\" \"134 (36.8)\"\nSet-RowValue $t \"QUADRIL\" \"84 (23.1)\"\nSet-RowValue $t \"TRAUMA\" \"89 (24.5)\"\n"}
